# Update frequency table values (rows 2-5, columns B-X) with new run results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0,0.00447093889716841,0.0439642324888227,0.00745156482861401,0.994038748137109,0.0476900149031297,0.0312965722801788,0.93591654247392,0.000745156482861401,0.0022354694485842,0.737704918032787,0,0.97317436661699,0.959761549925484,0.230998509687034,0.590163934426229,0.968703427719821,0.0253353204172876,0.0029806259314456,0.0134128166915052,0.168405365126677,0.00670640834575261,0.011177347242921),
    @(0.967213114754098,0.959761549925484,0.00745156482861401,0.0320417287630402,0.0029806259314456,0.011177347242921,0.000745156482861401,0,0.979135618479881,0.000745156482861401,0.0014903129657228,0.968703427719821,0.0029806259314456,0.0245901639344262,0.00819672131147541,0.0290611028315946,0.0268256333830104,0.931445603576751,0.011177347242921,0.000745156482861401,0.789865871833085,0.0320417287630402,0.0327868852459016),
    @(0.0014903129657228,0.0268256333830104,0.922503725782414,0.960506706408346,0.0014903129657228,0.927719821162444,0.0029806259314456,0.00521609538002981,0.000745156482861401,0,0.00819672131147541,0.00447093889716841,0.018628912071535,0.0119225037257824,0.760059612518629,0.3725782414307,0.003725782414307,0.00894187779433681,0.0253353204172876,0.982861400894188,0.0275707898658718,0.0208643815201192,0.938897168405365),
    @(0.0312965722801788,0.00894187779433681,0.026080476900149,0,0.0014903129657228,0.0134128166915052,0.964977645305514,0.0588673621460507,0.0193740685543964,0.997019374068554,0.252608047690015,0.0268256333830104,0.00521609538002981,0.003725782414307,0,0.00819672131147541,0.000745156482861401,0.033532041728763,0.960506706408346,0.0029806259314456,0.0141579731743666,0.940387481371088,0.0171385991058122)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2   # rows 2..5
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $colIndex = $j + 2   # columns B..X
        $ws.Cells.Item($rowIndex, $colIndex).Value = $rowValues[$j]
    }
}

$wb.Save()
